# Apply the "Done with 1873. Calculate Special Bonus" edit to ProbList.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C: enable wrap-text formatting (matches B and H columns) ---
$ws.Range("C1:C5").WrapText = $true

# --- New data row 6: "Calculate Special Bonus" problem (#1873) ---
$ws.Range("A6").Value = 1873
$ws.Range("B6").Value = "Calculate Special Bonus"
$ws.Range("C6").Value = "SELECT and ORDER"
$ws.Range("E6").Value = "Easy"
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = "✅"
$ws.Range("H6").Value = "Given 1 sol and didn’t see solutions"

$ws.Range("B6:C6").WrapText = $true
$ws.Range("H6").WrapText = $true

$ws.Rows.Item(6).RowHeight = 30

# --- Recolor the "done" highlight fill used by column A (green -> blue) ---
# (apply to A5:A6 first, then A2:A4, so the shared style slot is reused/mutated
#  in place instead of spawning extra, unreferenced style entries)
$ws.Range("A5:A6").Interior.Color = 15773696
$ws.Range("A2:A4").Interior.Color = 15773696

# --- Update the active selection to reflect where editing left off ---
$ws.Range("D10").Select() | Out-Null
